$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ09779051",
    "summ09981079",
    "summ10175210",
    "summ10372236",
    "summ10584754",
    "summ10802780",
    "summ11340166",
    "summ11709331",
    "summ12047680",
    "summ12355483",
    "summ12718705",
    "summ13008552",
    "summ13262608",
    "summ13558815",
    "summ13834316",
    "summ14093341",
    "summ14380891",
    "summ14640520",
    "summ14901562",
    "summ15143086",
    "summ15412430",
    "summ15683743",
    "summ15951351",
    "summ16244119",
    "summ16525019",
    "summ16785037",
    "summ17092010",
    "summ17378411",
    "summ17650590",
    "summ17908028",
    "summ18202380",
    "summ18464132",
    "summ18701792",
    "summ18952638",
    "summ19209193",
    "summ19508190",
    "summ19796391",
    "summ20057376",
    "summ20297499",
    "summ20577443",
    "summ20821484",
    "summ21073506",
    "summ21392030",
    "summ21647895",
    "summ21909437",
    "summ22206236",
    "summ22458399",
    "summ22721146",
    "summ23014133",
    "summ23328844"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

Write-Output "Renamed $($wb.Worksheets.Count) sheets"
